$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 34
$ws.Range("Q34").Value = 494
$ws.Range("Y34").Value = 1159

# Add new row 35 with data
# A35 holds a text string that looks like a date ("01-04-2021"). Assigning it
# directly via .Value would make Excel auto-convert it into a date serial
# number with a date number format. To store it as plain text (matching the
# other cells in column A), write it as a formula producing that text, then
# convert the formula to a static value via copy / paste-special-values -
# this avoids the date auto-detection while still ending with a plain text
# cell that uses the default style.
$ws.Range("A35").Formula = '="01-04-2021"'
$ws.Range("A35").Copy()
$ws.Range("A35").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("B35").Value = 343
$ws.Range("C35").Value = 45
$ws.Range("D35").Value = -541
$ws.Range("E35").Value = 65
$ws.Range("F35").Value = 354
$ws.Range("G35").Value = 32
$ws.Range("H35").Value = -24
$ws.Range("I35").Value = 33
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 161
$ws.Range("L35").Value = -110
$ws.Range("M35").Value = 78
$ws.Range("N35").Value = 954
$ws.Range("O35").Value = 220
$ws.Range("P35").Value = -752
$ws.Range("Q35").Value = -64
$ws.Range("R35").Value = 2486
$ws.Range("S35").Value = 114
$ws.Range("T35").Value = -51
$ws.Range("U35").Value = 293
$ws.Range("V35").Value = 4137
$ws.Range("W35").Value = 571
$ws.Range("X35").Value = -1478
$ws.Range("Y35").Value = 406
